{"js": "// Apply the three \"Guia exposici\u00f2 inicial\" text edits using the Word\n// JavaScript API (Office.js). `context` (alias `ctx`) is the\n// Word.RequestContext passed in by the harness.\n\nconst body = context.document.body;\n\n// --- Change 1 -----------------------------------------------------------\n// \"...la cual usaremos para hacer pruebas.\" ->\n// \"...la cual usaremos para hacer pruebas y comprobar el correcto\n//  funcionamiento de la instalaci\u00f3n.\"\nconst hit1 = body.search(\"para hacer pruebas.\", { matchCase: true });\nhit1.load(\"items\");\nawait context.sync();\nif (hit1.items.length === 0) {\n  throw new Error(\"No se encontr\u00f3 el texto de la tecnolog\u00eda Python/Pip.\");\n}\nhit1.items[0].insertText(\n  \"para hacer pruebas y comprobar el correcto funcionamiento de la instalaci\u00f3n.\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- Change 2 -------------------------------------------------------------\n// \"Visual Studio Code, un entorno que nos permite manejar tanto Docker como\n//  Python gracias a sus plugins.\" ->\n// \"Visual Studio Code, un IDE que nos permite manejar tanto Docker como\n//  Python 3 gracias a sus plugins.\"\nconst hit2 = body.search(\n  \"Visual Studio Code, un entorno que nos permite manejar tanto Docker como Python gracias a sus plugins.\",\n  { matchCase: true }\n);\nhit2.load(\"items\");\nawait context.sync();\nif (hit2.items.length === 0) {\n  throw new Error(\"No se encontr\u00f3 el p\u00e1rrafo de Visual Studio Code.\");\n}\nhit2.items[0].insertText(\n  \"Visual Studio Code, un IDE que nos permite manejar tanto Docker como Python 3 gracias a sus plugins.\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- Change 3 -------------------------------------------------------------\n// \"Finalmente, hemos pensado en Git, ya que nos ser\u00e1 \u00fatil a la hora de\n//  manejar ficheros y podr\u00eda ser la herramienta que use la empresa para\n//  gestionar sus proyectos.\" ->\n// \"Finalmente, hemos pensado en Git, una reconocida herramienta de control\n//  de versiones usada profesionalmente para gestionar proyectos. Adem\u00e1s,\n//  nuestra intenci\u00f3n ser\u00eda el uso de GIT para trasladar los proyectos del\n//  entorno Python 3 a Docker.\"\nconst hit3 = body.search(\n  \"ya que nos ser\u00e1 \u00fatil a la hora de manejar ficheros y podr\u00eda ser la herramienta que use la empresa para gestionar sus proyectos.\",\n  { matchCase: true }\n);\nhit3.load(\"items\");\nawait context.sync();\nif (hit3.items.length === 0) {\n  throw new Error(\"No se encontr\u00f3 el p\u00e1rrafo de Git.\");\n}\nhit3.items[0].insertText(\n  \"una reconocida herramienta de control de versiones usada profesionalmente para gestionar proyectos. Adem\u00e1s, nuestra intenci\u00f3n ser\u00eda el uso de GIT para trasladar los proyectos del entorno Python 3 a Docker.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Apply the three \"Guia exposici\u00f2 inicial\" text edits using the Word COM\n# object model. $word / $app resolve to the running application and the\n# document under edit is $word.ActiveDocument.\n\n$d = $word.ActiveDocument\n\n# --- Change 1 -------------------------------------------------------------\n# \"...la cual usaremos para hacer pruebas.\" ->\n# \"...la cual usaremos para hacer pruebas y comprobar el correcto\n#  funcionamiento de la instalaci\u00f3n.\"\n$r1 = $d.Content\n$r1.Find.Execute(\n    \"para hacer pruebas.\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"para hacer pruebas y comprobar el correcto funcionamiento de la instalaci\u00f3n.\",\n    2\n)\n\n# --- Change 2 ---------------------------------------------------------------\n# \"Visual Studio Code, un entorno que nos permite manejar tanto Docker como\n#  Python gracias a sus plugins.\" ->\n# \"Visual Studio Code, un IDE que nos permite manejar tanto Docker como\n#  Python 3 gracias a sus plugins.\"\n$r2 = $d.Content\n$r2.Find.Execute(\n    \"Visual Studio Code, un entorno que nos permite manejar tanto Docker como Python gracias a sus plugins.\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"Visual Studio Code, un IDE que nos permite manejar tanto Docker como Python 3 gracias a sus plugins.\",\n    2\n)\n\n# --- Change 3 ---------------------------------------------------------------\n# \"Finalmente, hemos pensado en Git, ya que nos ser\u00e1 \u00fatil a la hora de\n#  manejar ficheros y podr\u00eda ser la herramienta que use la empresa para\n#  gestionar sus proyectos.\" ->\n# \"Finalmente, hemos pensado en Git, una reconocida herramienta de control\n#  de versiones usada profesionalmente para gestionar proyectos. Adem\u00e1s,\n#  nuestra intenci\u00f3n ser\u00eda el uso de GIT para trasladar los proyectos del\n#  entorno Python 3 a Docker.\"\n$r3 = $d.Content\n$r3.Find.Execute(\n    \"ya que nos ser\u00e1 \u00fatil a la hora de manejar ficheros y podr\u00eda ser la herramienta que use la empresa para gestionar sus proyectos.\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"una reconocida herramienta de control de versiones usada profesionalmente para gestionar proyectos. Adem\u00e1s, nuestra intenci\u00f3n ser\u00eda el uso de GIT para trasladar los proyectos del entorno Python 3 a Docker.\",\n    2\n)\n"}
